$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MimsSoil")

# Reorder header row (row 1) columns A1:AR1 per target layout.
$ws.Range("A1").Value = "collection_date"
$ws.Range("B1").Value = "depth"
$ws.Range("C1").Value = "tax_class"
$ws.Range("D1").Value = "geo_loc_name"
$ws.Range("E1").Value = "project_name"
$ws.Range("F1").Value = "nucl_acid_ext"
$ws.Range("G1").Value = "lib_reads_seqd"
$ws.Range("H1").Value = "lat_lon"
$ws.Range("I1").Value = "env_local_scale"
$ws.Range("J1").Value = "samp_name"
$ws.Range("K1").Value = "elev"
$ws.Range("L1").Value = "sim_search_meth"
$ws.Range("M1").Value = "temp"
$ws.Range("N1").Value = "samp_taxon_id"
$ws.Range("O1").Value = "samp_mat_process"
$ws.Range("P1").Value = "lib_screen"
$ws.Range("Q1").Value = "seq_meth"
$ws.Range("R1").Value = "samp_size"
$ws.Range("S1").Value = "source_mat_id"
$ws.Range("T1").Value = "mid"
$ws.Range("U1").Value = "assembly_qual"
$ws.Range("V1").Value = "size_frac"
$ws.Range("W1").Value = "env_medium"
$ws.Range("X1").Value = "samp_collect_device"
$ws.Range("Y1").Value = "feat_pred"
$ws.Range("Z1").Value = "lib_size"
$ws.Range("AA1").Value = "env_broad_scale"
$ws.Range("AB1").Value = "lib_vector"
$ws.Range("AC1").Value = "assembly_name"
$ws.Range("AD1").Value = "samp_vol_we_dna_ext"
$ws.Range("AE1").Value = "adapters"
$ws.Range("AF1").Value = "number_contig"
$ws.Range("AG1").Value = "neg_cont_type"
$ws.Range("AH1").Value = "nucl_acid_amp"
$ws.Range("AI1").Value = "alt"
$ws.Range("AJ1").Value = "lib_layout"
$ws.Range("AK1").Value = "annot"
$ws.Range("AL1").Value = "experimental_factor"
$ws.Range("AM1").Value = "pos_cont_type"
$ws.Range("AN1").Value = "ref_biomaterial"
$ws.Range("AO1").Value = "assembly_software"
$ws.Range("AP1").Value = "rel_to_oxygen"
$ws.Range("AQ1").Value = "ref_db"
$ws.Range("AR1").Value = "samp_collect_method"

# Move data validations to follow their associated (now relocated) columns.
function Move-Validation($ws, $oldSqref, $newSqref) {
    $oldRange = $ws.Range($oldSqref)
    $validation = $oldRange.Validation
    $vType = $validation.Type()
    $vFormula1 = $validation.Formula1()
    $validation.Delete()
    $newRange = $ws.Range($newSqref)
    $newValidation = $newRange.Validation
    $newValidation.Add($vType, 1, 1, '"' + $vFormula1 + '"')
    $newValidation.IgnoreBlank = $true
    $newValidation.InCellDropdown = $true
    $newValidation.ShowInput = $false
    $newValidation.ShowError = $false
}

Move-Validation $ws "N2:N1048576" "AG2:AG1048576"
Move-Validation $ws "Y2:Y1048576" "AJ2:AJ1048576"
Move-Validation $ws "AH2:AH1048576" "AP2:AP1048576"

Write-Host "Header reorder and data validation relocation complete."
